# functional_composition_of_core_genes.xlsx - "Subsystem" sheet refresh
# Commit: "results with fixed workflow" - re-run of the classification workflow
# produced 4 additional subsystem categories (Biomass and maintenance functions,
# Extracellular exchange, Intracellular demand, Intracellular source/sink) and
# recomputed every BAR/NBR percentage against the new gene-set totals, so every
# data row (2-56) is rewritten in place with its final label/value triple.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subsystem")
$ws.Activate()

$ws.Cells.Item(2, 1).Value = "Alanine, aspartate and glutamate metabolism"
$ws.Cells.Item(2, 2).Value = 1.482701812191104
$ws.Cells.Item(2, 3).Value = 2.325581395348837
$ws.Cells.Item(3, 1).Value = "Aminosugars metabolism"
$ws.Cells.Item(3, 2).Value = 0.8237232289950577
$ws.Cells.Item(3, 3).Value = 0.7751937984496124
$ws.Cells.Item(4, 1).Value = "Arginine and proline metabolism"
$ws.Cells.Item(4, 2).Value = 2.635914332784185
$ws.Cells.Item(4, 3).Value = 10.85271317829457
$ws.Cells.Item(5, 1).Value = "Biomass and maintenance functions"
$ws.Cells.Item(5, 2).Value = 0.1647446457990115
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 1).Value = "Biotin metabolism"
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 2.325581395348837
$ws.Cells.Item(7, 1).Value = "C5-Branched dibasic acid metabolism"
$ws.Cells.Item(7, 2).Value = 0.4942339373970346
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 1).Value = "Carbon fixation"
$ws.Cells.Item(8, 2).Value = 0.9884678747940692
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 1).Value = "Carotenoid Biosynthesis"
$ws.Cells.Item(9, 2).Value = 2.14168039538715
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 1).Value = "Citrate cycle (TCA cycle)"
$ws.Cells.Item(10, 2).Value = 0.9884678747940692
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 1).Value = "Cyanophycin metabolism"
$ws.Cells.Item(11, 2).Value = 0.3294892915980231
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(12, 1).Value = "Extracellular exchange"
$ws.Cells.Item(12, 2).Value = 3.130148270181219
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 1).Value = "Fatty acid biosynthesis"
$ws.Cells.Item(13, 2).Value = 16.80395387149918
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 1).Value = "Folate biosynthesis"
$ws.Cells.Item(14, 2).Value = 2.14168039538715
$ws.Cells.Item(14, 3).Value = 1.550387596899225
$ws.Cells.Item(15, 1).Value = "Fructose and mannose metabolism"
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 6.2015503875969
$ws.Cells.Item(16, 1).Value = "Galactolipids metabolism"
$ws.Cells.Item(16, 2).Value = 4.448105436573312
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(17, 1).Value = "Glutamate metabolism"
$ws.Cells.Item(17, 2).Value = 1.153212520593081
$ws.Cells.Item(17, 3).Value = 0.7751937984496124
$ws.Cells.Item(18, 1).Value = "Glutathione metabolism"
$ws.Cells.Item(18, 2).Value = 0.1647446457990115
$ws.Cells.Item(18, 3).Value = 1.550387596899225
$ws.Cells.Item(19, 1).Value = "Glycerolipid metabolism"
$ws.Cells.Item(19, 2).Value = 0.6589785831960462
$ws.Cells.Item(19, 3).Value = 1.550387596899225
$ws.Cells.Item(20, 1).Value = "Glycolysis/Gluconeogenesis"
$ws.Cells.Item(20, 2).Value = 3.130148270181219
$ws.Cells.Item(20, 3).Value = 4.651162790697675
$ws.Cells.Item(21, 1).Value = "Glyoxylate and dicarboxylate metabolism"
$ws.Cells.Item(21, 2).Value = 2.471169686985173
$ws.Cells.Item(21, 3).Value = 0.7751937984496124
$ws.Cells.Item(22, 1).Value = "Histidine metabolism"
$ws.Cells.Item(22, 2).Value = 1.482701812191104
$ws.Cells.Item(22, 3).Value = 1.550387596899225
$ws.Cells.Item(23, 1).Value = "Hydrogen production"
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 0.7751937984496124
$ws.Cells.Item(24, 1).Value = "Inositol phosphate metabolism"
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 1.550387596899225
$ws.Cells.Item(25, 1).Value = "Intracellular demand"
$ws.Cells.Item(25, 2).Value = 0.1647446457990115
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(26, 1).Value = "Intracellular source/sink"
$ws.Cells.Item(26, 2).Value = 0.1647446457990115
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(27, 1).Value = "Lipopolysaccharide biosynthesis"
$ws.Cells.Item(27, 2).Value = 0.8237232289950577
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(28, 1).Value = "Lysine metabolism"
$ws.Cells.Item(28, 2).Value = 1.482701812191104
$ws.Cells.Item(28, 3).Value = 1.550387596899225
$ws.Cells.Item(29, 1).Value = "Nicotinate and nicotinamide metabolism"
$ws.Cells.Item(29, 2).Value = 0.8237232289950577
$ws.Cells.Item(29, 3).Value = 1.550387596899225
$ws.Cells.Item(30, 1).Value = "Nitrogen metabolism"
$ws.Cells.Item(30, 2).Value = 1.812191103789127
$ws.Cells.Item(30, 3).Value = 0.7751937984496124
$ws.Cells.Item(31, 1).Value = "Nucleotide sugars metabolism"
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 3.10077519379845
$ws.Cells.Item(32, 1).Value = "Others"
$ws.Cells.Item(32, 2).Value = 0.6589785831960462
$ws.Cells.Item(32, 3).Value = 1.550387596899225
$ws.Cells.Item(33, 1).Value = "Oxidative phosphorylation"
$ws.Cells.Item(33, 2).Value = 0.8237232289950577
$ws.Cells.Item(33, 3).Value = 1.550387596899225
$ws.Cells.Item(34, 1).Value = "PHB byosynthesis"
$ws.Cells.Item(34, 2).Value = 0
$ws.Cells.Item(34, 3).Value = 2.325581395348837
$ws.Cells.Item(35, 1).Value = "Pantothenate and CoA biosynthesis"
$ws.Cells.Item(35, 2).Value = 1.482701812191104
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(36, 1).Value = "Pentose phosphate pathway"
$ws.Cells.Item(36, 2).Value = 1.482701812191104
$ws.Cells.Item(36, 3).Value = 0.7751937984496124
$ws.Cells.Item(37, 1).Value = "Peptidoglycan biosynthesis"
$ws.Cells.Item(37, 2).Value = 1.317957166392092
$ws.Cells.Item(37, 3).Value = 0.7751937984496124
$ws.Cells.Item(38, 1).Value = "Phenylalanine tyrosine and tryptophan biosynthesis"
$ws.Cells.Item(38, 2).Value = 3.130148270181219
$ws.Cells.Item(38, 3).Value = 5.426356589147287
$ws.Cells.Item(39, 1).Value = "Photosynthesis"
$ws.Cells.Item(39, 2).Value = 1.482701812191104
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(40, 1).Value = "Porphyrin and chlorophyll metabolism"
$ws.Cells.Item(40, 2).Value = 6.754530477759473
$ws.Cells.Item(40, 3).Value = 5.426356589147287
$ws.Cells.Item(41, 1).Value = "Purine metabolism"
$ws.Cells.Item(41, 2).Value = 4.612850082372323
$ws.Cells.Item(41, 3).Value = 6.2015503875969
$ws.Cells.Item(42, 1).Value = "Pyrimidine metabolism"
$ws.Cells.Item(42, 2).Value = 3.789126853377265
$ws.Cells.Item(42, 3).Value = 2.325581395348837
$ws.Cells.Item(43, 1).Value = "Pyruvate metabolism"
$ws.Cells.Item(43, 2).Value = 2.14168039538715
$ws.Cells.Item(43, 3).Value = 0.7751937984496124
$ws.Cells.Item(44, 1).Value = "Riboflavin metabolism"
$ws.Cells.Item(44, 2).Value = 1.647446457990115
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(45, 1).Value = "Starch and sucrose metabolism"
$ws.Cells.Item(45, 2).Value = 0.9884678747940692
$ws.Cells.Item(45, 3).Value = 3.875968992248062
$ws.Cells.Item(46, 1).Value = "Steroid biosynthesis"
$ws.Cells.Item(46, 2).Value = 0
$ws.Cells.Item(46, 3).Value = 2.325581395348837
$ws.Cells.Item(47, 1).Value = "Sterol biosynthesis"
$ws.Cells.Item(47, 2).Value = 0.9884678747940692
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(48, 1).Value = "Sulfolipid Biosynthesis"
$ws.Cells.Item(48, 2).Value = 1.647446457990115
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(49, 1).Value = "Sulfur Cysteine and methionine metabolism"
$ws.Cells.Item(49, 2).Value = 2.306425041186162
$ws.Cells.Item(49, 3).Value = 3.875968992248062
$ws.Cells.Item(50, 1).Value = "Terpenoid backbone biosynthesis"
$ws.Cells.Item(50, 2).Value = 1.482701812191104
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(51, 1).Value = "Thiamine metabolism"
$ws.Cells.Item(51, 2).Value = 0
$ws.Cells.Item(51, 3).Value = 3.10077519379845
$ws.Cells.Item(52, 1).Value = "Transport"
$ws.Cells.Item(52, 2).Value = 6.919275123558484
$ws.Cells.Item(52, 3).Value = 7.751937984496124
$ws.Cells.Item(53, 1).Value = "Ubiquinone and other pterpenoids biosynthesis"
$ws.Cells.Item(53, 2).Value = 2.306425041186162
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(54, 1).Value = "Urea cycle and metabolism of amino groups"
$ws.Cells.Item(54, 2).Value = 1.153212520593081
$ws.Cells.Item(54, 3).Value = 0.7751937984496124
$ws.Cells.Item(55, 1).Value = "Valine leucine and isoleucine biosynthesis"
$ws.Cells.Item(55, 2).Value = 1.976935749588138
$ws.Cells.Item(55, 3).Value = 3.10077519379845
$ws.Cells.Item(56, 1).Value = "Vitamin B6 metabolism"
$ws.Cells.Item(56, 2).Value = 0
$ws.Cells.Item(56, 3).Value = 3.875968992248062

# The four newly-inserted rows fall outside the workbook's original used range,
# so their column-A cells start out unstyled. Copy the header-row formatting
# (bold font, thin border, centered alignment) from an existing labeled cell
# so the new rows match the rest of column A.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A53:A56").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select()
